# "can select skip to generate any enum"
#
# The sample data sheet ("Sheet1", the 3rd tab) gets a second copy of the
# demo table (currently A1:O7) pasted below it at A11:O17, so the sample
# shows that the enum code-gen flag can be toggled per-column instead of
# only applying it uniformly across the whole row: the new header row
# keeps the `{true}` (skip) marker on the "e0" column but drops it again
# for "e1"/"e2", reusing the plain `TestEnum:e1` / `TestEnum:e2` header
# strings that already exist elsewhere in the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteValues  = [Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues
$xlPasteFormats = [Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats

# Duplicate the whole demo table (header + 5 sample rows) ten rows down.
$ws.Range("A1:O7").Copy()
$ws.Range("A11").PasteSpecial($xlPasteValues)

# The source table's "is_boss" column (D) uses a quote-prefixed cell
# style (so TRUE/FALSE literals don't get auto-typed); carry that format
# down to the pasted copy as well.
$ws.Range("D3:D7").Copy()
$ws.Range("D13:D17").PasteSpecial($xlPasteFormats)

# On the new header row, show that individual enum columns can opt back
# out of the "{true}" skip flag while a sibling column (M) keeps it.
$ws.Range("N12").Value2 = "TestEnum:e1"
$ws.Range("O12").Value2 = "TestEnum:e2"

# Leave the new block selected, like after pasting it in.
$ws.Range("A11:O17").Select()
